$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Defined names: "loylecarner1" (Sheet1/Sheet3 scoped) -> "mavisstaples4",
#    and shrink the referenced range from $A$1:$E$15 to $A$1:$E$12 (new web
#    query table has fewer rows). Renaming a Name object in-place while a
#    same-named sibling still exists is unreliable here, so delete + re-add.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

$names = $wb.Names
while ($names.Count -gt 0) {
  $names.Item(1).Delete()
}

$ws1.Names.Add("mavisstaples4", "=Sheet1!`$A`$1:`$E`$12")
$ws3.Names.Add("mavisstaples4", "=Sheet3!`$A`$1:`$E`$12")

# ---------------------------------------------------------------------------
# 2. Replace the web-query result data on Sheet1 & Sheet3 (identical data):
#    old "Loyle Carner" tracklist (14 rows) -> new "Mavis Staples" tracklist
#    (10 rows), refreshed from the renamed query/connection.
# ---------------------------------------------------------------------------
foreach ($sheetName in @("Sheet1", "Sheet3")) {
  $ws = $wb.Worksheets.Item($sheetName)

  $ws.Range("B2").Value = "Little Bit"
  $ws.Range("C2").Value = "Jeff Tweedy"
  $ws.Range("D2").Value = "Mavis Staples"
  $ws.Range("E2").Value = 0.16041666666666668

  $ws.Range("B3").Value = "If All I Was Was Black"
  $ws.Range("C3").Value = "Mavis Staples / Jeff Tweedy"
  $ws.Range("D3").Value = "Mavis Staples"
  $ws.Range("E3").Value = 0.16319444444444445

  $ws.Range("B4").Value = "Who Told You That"
  $ws.Range("C4").Value = "Jeff Tweedy"
  $ws.Range("D4").Value = "Mavis Staples"
  $ws.Range("E4").Value = 0.11666666666666665

  $ws.Range("B5").Value = "Ain't No Doubt About It"
  $ws.Range("C5").Value = "Jeff Tweedy"
  $ws.Range("D5").Value = "Mavis Staples"
  $ws.Range("E5").Value = 0.13749999999999998

  $ws.Range("B6").Value = "Peaceful Dream"
  $ws.Range("C6").Value = "Jeff Tweedy"
  $ws.Range("D6").Value = "Mavis Staples"
  $ws.Range("E6").Value = 0.1388888888888889

  $ws.Range("B7").Value = "No Time for Crying"
  $ws.Range("C7").Value = "Mavis Staples / Jeff Tweedy"
  $ws.Range("D7").Value = "Mavis Staples"
  $ws.Range("E7").Value = 0.19166666666666665

  $ws.Range("B8").Value = "Build a Bridge"
  $ws.Range("C8").Value = "Jeff Tweedy"
  $ws.Range("D8").Value = "Mavis Staples"
  $ws.Range("E8").Value = 0.15069444444444444

  $ws.Range("B9").Value = "We Go High"
  $ws.Range("C9").Value = "Mavis Staples / Jeff Tweedy"
  $ws.Range("D9").Value = "Mavis Staples"
  $ws.Range("E9").Value = 0.14305555555555557

  $ws.Range("B10").Value = "Try Harder"
  $ws.Range("C10").Value = "Jeff Tweedy"
  $ws.Range("D10").Value = "Mavis Staples"
  $ws.Range("E10").Value = 0.16041666666666668

  $ws.Range("B11").Value = "All Over Again"
  $ws.Range("C11").Value = "Jeff Tweedy"
  $ws.Range("D11").Value = "Mavis Staples"
  $ws.Range("E11").Value = 0.079166666666666663

  # Rows 12-15 held the remaining old tracks; the refreshed (smaller) table
  # no longer fills them, so the query clears their contents but leaves the
  # row/cell formatting (style) intact - ClearContents matches that.
  $ws.Range("A12:E15").ClearContents()

  # The query result area grew by a few trailing formatted (but empty) rows
  # at the bottom of the sheet (G column keeps its time-style formatting).
  $ws.Range("G29:G33").NumberFormat = "h:mm"

  # Column widths re-auto-fit around the new (narrower) data.
  $ws.Columns.Item(2).ColumnWidth = 20.9331
  $ws.Columns.Item(3).ColumnWidth = 25.6513
  $ws.Columns.Item(4).ColumnWidth = 12.6499
  $ws.Columns.Item(6).ColumnWidth = 25.9374
}

# ---------------------------------------------------------------------------
# 3. Sheet2 is a formula-driven report sheet (LEFTB/REPT/MAX(LENB(...)))
#    that reads straight off Sheet1, so it recalculates on its own. Only its
#    UI selection needs to move to reflect the new (shorter) data range.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate() | Out-Null
$ws2.Range("A3:K14").Select() | Out-Null
